$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.090.91'
$ws.Range('D3').Value = '2.938.70'
$ws.Range('E3').Value = '  +4.50%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '203.03'
$ws.Range('E5').Value = '  +8.13%  '
$ws.Range('D6').Value = '599.80'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('E9').Value = '  +5.05%  '
$ws.Range('D10').Value = '2.940.62'
$ws.Range('E10').Value = '  +4.62%  '
$ws.Range('D11').Value = '0.449'
$ws.Range('E11').Value = '  +19.68%  '
$ws.Range('D12').Value = '0.162'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = '3.487.25'
$ws.Range('D15').Value = '28.49'
$ws.Range('E15').Value = '  +6.21%  '
$ws.Range('D16').Value = '76.031.16'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('D18').Value = '2.937.36'
$ws.Range('E18').Value = '  +4.35%  '
$ws.Range('D19').Value = '13.33'
$ws.Range('E19').Value = '  +8.84%  '
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').Value = '375.35'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').Value = '4.37'
$ws.Range('E23').Value = '  +7.17%  '
$ws.Range('D24').Value = '71.88'
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '4.38'
$ws.Range('E26').Value = '  +5.78%  '
$ws.Range('D27').Value = '3.094.96'
$ws.Range('E27').Value = '  +4.84%  '
$ws.Range('D28').Value = '9.78'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('E29').Value = '  +7.20%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '7.97'
$ws.Range('E32').Value = '  +5.09%  '
$ws.Range('D33').Value = '503.64'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('E34').Value = '  +4.11%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '165.18'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '20.39'
$ws.Range('E37').Value = '  +2.87%  '
$ws.Range('D38').Value = '0.110'
$ws.Range('E38').Value = '  +27.24%  '
$ws.Range('D39').Value = '19.66'
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('D40').Value = '0.375'
$ws.Range('E40').Value = '  +10.36%  '
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').Value = '180.66'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '5.04'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = '1.21'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '40.14'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('E48').Value = '  +2.36%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '3.89'
$ws.Range('E49').Value = '  +4.99%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '0.585'
$ws.Range('E50').Value = '  +2.32%  '
$ws.Range('D51').Value = '22.84'
$ws.Range('E51').Value = '  +9.96%  '
